$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 126, pushing the former rows 126 and 127
# (the "Evan Sooklal vs Lilith Karyadi" and "Lilith Karyadi vs Will Simpson"
# match results) down to rows 128 and 129 respectively.
$ws.Rows(126).Insert()
$ws.Rows(126).Insert()

# The two freshly inserted rows come back with generic/default formatting.
# Copy the column-A cell format (border/alignment/font) from row 125 so the
# new index cells match the rest of the table.
$ws.Range("A125").Copy()
$ws.Range("A126:A127").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 126: new match result between Roman Ramirez (G) and Will Simpson (H).
$ws.Range("A126").Value = 125
$ws.Range("G126").Value = 1185.628639295045
$ws.Range("H126").Value = 1222.322874151305

# Row 127: new match result between Roman Ramirez (G) and Will Simpson (H).
$ws.Range("A127").Value = 126
$ws.Range("G127").Value = 1173.211242489557
$ws.Range("H127").Value = 1234.740270956793

# Row 128 (formerly row 126) keeps its Evan Sooklal (C) / Lilith Karyadi (E)
# result, only its running index in column A needs to move from 125 to 127.
$ws.Range("A128").Value = 127

# Row 129 (formerly row 127) keeps the Lilith Karyadi (E) / Will Simpson (H)
# pairing, but the Elo values are recalculated because of the two newly
# inserted matches earlier in the table.
$ws.Range("A129").Value = 128
$ws.Range("E129").Value = 1228.001242930469
$ws.Range("H129").Value = 1251.05480171003
